$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Date_of_sales (column I) for all data rows 2-34: 45371 -> 45417
$ws.Range("I2:I34").Value = 45417

# Reorder words/tokens in columns C and G so that "сер" / "сер," comes first
$ws.Range("C2").Value = "сер легк б/к"
$ws.Range("G2").Value = "сер, легк, б/к"
$ws.Range("C3").Value = "сер легк б/к"
$ws.Range("G3").Value = "сер, легк, б/к"
$ws.Range("C5").Value = "сер легк"
$ws.Range("G5").Value = "сер, легк"
$ws.Range("C6").Value = "сер легк"
$ws.Range("G6").Value = "сер, легк"
$ws.Range("C7").Value = "сер легк"
$ws.Range("G7").Value = "сер, легк"
$ws.Range("C8").Value = "сер легк"
$ws.Range("G8").Value = "сер, легк"
$ws.Range("C9").Value = "210B сер Type H"
$ws.Range("G9").Value = "сер, легк"
$ws.Range("C10").Value = "сер б/к груз"
$ws.Range("G10").Value = "210B, сер, Type, H"
$ws.Range("C11").Value = "сер легк б/к"
$ws.Range("G11").Value = "сер, б/к, груз"
$ws.Range("C12").Value = "сер легк б/к"
$ws.Range("G12").Value = "сер, б/к, груз"
$ws.Range("C13").Value = "8 сер сх"
$ws.Range("G13").Value = "сер, б/к, груз"
$ws.Range("C14").Value = "сер легк"
$ws.Range("G14").Value = "сер, б/к, груз"
$ws.Range("C15").Value = "сер легк"
$ws.Range("G15").Value = "сер, легк, б/к"
$ws.Range("C16").Value = "сер легк"
$ws.Range("G16").Value = "сер, легк, б/к"
$ws.Range("C17").Value = "сер легк"
$ws.Range("G17").Value = "сер, легк"
$ws.Range("C18").Value = "сер легк"
$ws.Range("G18").Value = "сер, легк"
$ws.Range("C19").Value = "сер легк"
$ws.Range("G19").Value = "сер, легк"
$ws.Range("C20").Value = "сер легк б/к"
$ws.Range("G20").Value = "сер, легк"
$ws.Range("C21").Value = "сер легк б/к"
$ws.Range("G21").Value = "сер, легк"
$ws.Range("C22").Value = "сер легк б/к"
$ws.Range("G22").Value = "сер, легк"
$ws.Range("G23").Value = "сер, легк"
$ws.Range("C24").Value = "сер легк б/к"
$ws.Range("G24").Value = "сер, легк"
$ws.Range("C25").Value = "сер легк"
$ws.Range("G25").Value = "сер, легк, б/к"
$ws.Range("G26").Value = "сер, легк, б/к"
$ws.Range("G27").Value = "сер, легк, б/к"
$ws.Range("G29").Value = "сер, легк, б/к"
$ws.Range("G30").Value = "сер, легк, б/к"
$ws.Range("G31").Value = "сер, легк, б/к"
$ws.Range("G32").Value = "сер, легк, б/к"
$ws.Range("G33").Value = "сер, легк, б/к"
$ws.Range("G34").Value = "сер, легк"
